# Updated footprints to 0805 packages
#
# Rows 5-9 (R1/R2/R7/R8/R9 = 2.2K, R3 = 10K, R4 = 1K, R5/R6 = 220R, C1 = 0.1uF)
# get new Digikey part numbers / manufacturer part numbers / descriptions /
# datasheet links reflecting a switch from 0603 to 0805 footprints.
# Everything else (quantities, unit prices, totals, other rows) is untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 5: R1,R2,R7,R8,R9 - 2.2K ---
$ws.Range("E5").Value = "311-2.20KCRCT-ND "
$ws.Range("F5").Value = "RC0805FR-072K2L"
$ws.Range("G5").Value = "RES 2.2K OHM 1% 1/8W 0805"

# --- Row 6: R3 - 10K ---
$ws.Range("E6").Value = "311-10.0KCRCT-ND"
$ws.Range("F6").Value = "RC0805FR-0710KL"
$ws.Range("G6").Value = "RES 10K OHM 1% 1/8W 0805"

# --- Row 7: R4 - 1K ---
$ws.Range("E7").Value = "311-1.00KCRCT-ND"
$ws.Range("F7").Value = "RC0805FR-071KL"
$ws.Range("G7").Value = "RES 1K OHM 1% 1/8W 0805"

# --- Row 8: R5,R6 - 220 ---
$ws.Range("E8").Value = "311-220ARCT-ND"
$ws.Range("F8").Value = "RC0805JR-07220RL"
$ws.Range("G8").Value = "RES 220 OHM 5% 1/8W 0805"

# --- Row 9: C1 - 0.1uF ---
$ws.Range("E9").Value = "478-3755-1-ND"
$ws.Range("F9").Value = "08053C104KAT2A"
$ws.Range("G9").Value = "CAP CER 0.1UF 25V X7R 0805"

# --- Update hyperlinks on column J for the changed rows (new target URLs) ---
# Drop every existing hyperlink and re-create them so both the changed ones
# (J5-J9) and the untouched ones (J2,J3,J4,J10,J11) come back with correct
# targets, in the order Excel ends up with after re-saving.
$ws.Range("A1").Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("J8"), "https://www.digikey.ca/en/products/detail/yageo/RC0805JR-07220RL/728283")
$ws.Hyperlinks.Add($ws.Range("J10"), "https://www.digikey.ca/en/products/detail/lite-on-inc/LTST-C190KRKT/386817")
$ws.Hyperlinks.Add($ws.Range("J11"), "https://www.digikey.ca/en/products/detail/omron-electronics-inc-emc-div/G3VM-41DY1-TR05/5799757 ")
$ws.Hyperlinks.Add($ws.Range("J2"), "https://www.digikey.ca/en/products/detail/cui-devices/SJ1-3513/738683 ")
$ws.Hyperlinks.Add($ws.Range("J3"), "https://www.digikey.ca/en/products/detail/texas-instruments/TCA9534DWR/6566100 ")
$ws.Hyperlinks.Add($ws.Range("J4"), "https://www.digikey.ca/en/products/detail/sparkfun-electronics/PRT-14417/7652746 ")
$ws.Hyperlinks.Add($ws.Range("J7"), "https://www.digikey.ca/en/products/detail/yageo/RC0805FR-071KL/727444")
$ws.Hyperlinks.Add($ws.Range("J6"), "https://www.digikey.ca/en/products/detail/yageo/RC0805FR-0710KL/727535")
$ws.Hyperlinks.Add($ws.Range("J5"), "https://www.digikey.ca/en/products/detail/yageo/RC0805FR-072K2L/727676")
$ws.Hyperlinks.Add($ws.Range("J9"), "https://www.digikey.ca/en/products/detail/avx-corporation/08053C104KAT2A/1116281")

# --- Widen column J (Link) to fit the new, longer descriptions/links ---
$ws.Columns.Item(10).ColumnWidth = 92.7
